# C5-PowerPoint.pptx edit
#
# 1) Slide 6's table switches from the custom "Table_0" style to the
#    built-in "Medium Style 2 - Accent 1" table style
#    ({B35DAE46-E88C-4074-8EFD-AB764B7CE0FE}).
# 2) The deck's applied design theme changes from "Integral" to the
#    default "Office Theme" color palette (dk1/lt1/dk2/lt2/accent1-6/
#    hlink/folHlink) on the slide master's theme.

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 6 (shape 2 is the Google-Slides-imported table) ---
$tableSlide = $p.Slides.Item(6)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{B35DAE46-E88C-4074-8EFD-AB764B7CE0FE}")

# --- 2) Swap the theme color scheme from "Integral" to "Office Theme" ---
# Office theme RGB values (standard hex), converted below to the BGR
# integer encoding PowerPoint's ColorFormat.RGB expects.
$officeColors = @(
    0x000000,  # 1  dk1
    0xFFFFFF,  # 2  lt1
    0x44546A,  # 3  dk2
    0xE7E6E6,  # 4  lt2
    0x5B9BD5,  # 5  accent1
    0xED7D31,  # 6  accent2
    0xA5A5A5,  # 7  accent3
    0xFFC000,  # 8  accent4
    0x4472C4,  # 9  accent5
    0x70AD47,  # 10 accent6
    0x0563C1,  # 11 hlink
    0x954F72   # 12 folHlink
)

$themeColorScheme = $p.SlideMaster.Theme.ThemeColorScheme

for ($i = 1; $i -le $themeColorScheme.Count; $i++) {
    $hex = $officeColors[$i - 1]
    $r = [math]::Floor($hex / 0x10000) -band 0xFF
    $g = [math]::Floor($hex / 0x100) -band 0xFF
    $b = $hex -band 0xFF
    $bgr = ($b * 0x10000) + ($g * 0x100) + $r
    $themeColorScheme.Colors($i).RGB = $bgr
}
